$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.592.30"
$ws.Range("E2").Value = "  -1.01%  "
$ws.Range("D3").Value = "1.597.09"
$ws.Range("E3").Value = "  -2.06%  "
$ws.Range("E4").Value = "  +0.50%  "
$ws.Range("D5").Value = "'208.35"
$ws.Range("E5").Value = "  -1.29%  "
$ws.Range("D6").Value = "'0.504"
$ws.Range("E6").Value = "  -3.38%  "
$ws.Range("E7").Value = "  +0.54%  "
$ws.Range("D8").Value = "'22.31"
$ws.Range("E8").Value = "  -4.41%  "
$ws.Range("E10").Value = "  -3.25%  "
$ws.Range("D11").Value = "'0.0865"
$ws.Range("E11").Value = "  -1.77%  "
$ws.Range("D12").Value = "1.824.16"
$ws.Range("E12").Value = "  -2.01%  "
$ws.Range("D13").Value = "1.595.84"
$ws.Range("E13").Value = "  -2.07%  "
$ws.Range("E14").Value = "  -3.95%  "
$ws.Range("E15").Value = "  -4.48%  "
$ws.Range("D16").Value = "'63.47"
$ws.Range("E16").Value = "  -2.77%  "
$ws.Range("D17").Value = "27.597.49"
$ws.Range("E17").Value = "  -1.03%  "
$ws.Range("D18").Value = "'217.72"
$ws.Range("E18").Value = "  -5.19%  "
$ws.Range("D19").Value = "'7.38"
$ws.Range("E19").Value = "  -3.61%  "
$ws.Range("E20").Value = "  -3.64%  "
$ws.Range("E21").Value = "  +0.58%  "
$ws.Range("D22").Value = "'4.18"
$ws.Range("E22").Value = "  -3.50%  "
$ws.Range("D23").Value = "'9.65"
$ws.Range("E23").Value = "  -4.56%  "
$ws.Range("E24").Value = "  -2.34%  "
$ws.Range("D25").Value = "'153.30"
$ws.Range("E26").Value = "  +0.53%  "
$ws.Range("D27").Value = "'6.74"
$ws.Range("E27").Value = "  -2.36%  "
$ws.Range("D28").Value = "'15.08"
$ws.Range("E29").Value = "  -3.73%  "
$ws.Range("E30").Value = "  -1.29%  "
$ws.Range("E31").Value = "  -2.93%  "
$ws.Range("E32").Value = "  -4.28%  "
$ws.Range("D33").Value = "1.368.84"
$ws.Range("E33").Value = "  -1.75%  "
$ws.Range("D34").Value = "'2.96"
$ws.Range("E34").Value = "  -4.92%  "
$ws.Range("E35").Value = "  -3.52%  "
$ws.Range("E36").Value = "  -3.80%  "
$ws.Range("E38").Value = "  -2.60%  "
$ws.Range("D39").Value = "'0.540"
$ws.Range("E39").Value = "  -3.05%  "
$ws.Range("D40").Value = "'0.813"
$ws.Range("E40").Value = "  -4.50%  "
$ws.Range("E41").Value = "  +0.56%  "
$ws.Range("D42").Value = "'0.977"
$ws.Range("E42").Value = "  -3.46%  "
$ws.Range("E43").Value = "  -1.18%  "
$ws.Range("D44").Value = "'1.79"
$ws.Range("E44").Value = "  -3.59%  "
$ws.Range("D45").Value = "'64.09"
$ws.Range("E45").Value = "  -2.55%  "
$ws.Range("D46").Value = "1.734.36"
$ws.Range("E46").Value = "  -2.08%  "
$ws.Range("E47").Value = "  -1.55%  "
$ws.Range("D48").Value = "'87.90"
$ws.Range("E48").Value = "  -0.58%  "
$ws.Range("E49").Value = "  -3.22%  "
$ws.Range("D50").Value = "'0.0971"
$ws.Range("E50").Value = "  -4.44%  "
$ws.Range("D51").Value = "'0.0498"
$ws.Range("E51").Value = "  -1.03%  "
